$d = $word.ActiveDocument

# The title paragraph currently holds two runs (identical formatting):
#   Run A: "Week 6-12"
#   Run B: ".React JS-HOL"
# immediately followed by a bookmarkStart/bookmarkEnd pair named "_GoBack"
# (sitting right after Run B).
#
# Target layout (still two runs, same formatting), but the text/split point
# changes and the bookmark moves to sit between the two runs instead of
# after them:
#   Run A: "Week 7"
#   bookmarkStart/bookmarkEnd "_GoBack"   <-- moved here
#   Run B: "-12.React JS-HOL"

# Locate the first run's text precisely (rather than hard-coding offsets).
$findRange1 = $d.Content
$findRange1.Find.ClearFormatting()
[void]$findRange1.Find.Execute("Week 6-12", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)

# Replace "Week 6-12" with "Week 7" in place.
$findRange1.Text = "Week 7"
$splitPoint = $findRange1.End

# Locate the second run's text, now immediately following the edited text.
$findRange2 = $d.Range($splitPoint, $d.Content.End)
$findRange2.Find.ClearFormatting()
[void]$findRange2.Find.Execute(".React JS-HOL", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)

# Replace ".React JS-HOL" with "-12.React JS-HOL" in place.
$findRange2.Text = "-12.React JS-HOL"

# Move the "_GoBack" bookmark so it sits between the two runs (right after
# "Week 7", at the split point) instead of after the second run. Re-adding
# a bookmark named "_GoBack" repositions the single reserved "_GoBack"
# bookmark Word keeps.
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
